$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.581.31"
$ws.Range("E2").Value = "'  -4.69%  "

$ws.Range("D3").Value = "'3.265.92"
$ws.Range("E3").Value = "'  -4.92%  "

$ws.Range("E4").Value = "'  -0.03%  "

$ws.Range("D5").Value = "'554.45"
$ws.Range("E5").Value = "'  -3.23%  "

$ws.Range("D6").Value = "'184.31"
$ws.Range("E6").Value = "'  -2.22%  "

$ws.Range("E7").Value = "'  -0.13%  "

$ws.Range("E8").Value = "'  -2.40%  "

$ws.Range("D9").Value = "'3.260.94"
$ws.Range("E9").Value = "'  -4.78%  "

$ws.Range("E10").Value = "'  -7.44%  "

$ws.Range("E11").Value = "'  -4.22%  "

$ws.Range("D12").Value = "'47.22"
$ws.Range("E12").Value = "'  -6.87%  "

$ws.Range("E13").Value = "'  -5.50%  "

$ws.Range("E14").Value = "'  -4.86%  "

$ws.Range("D15").Value = "'624.12"
$ws.Range("E15").Value = "'  -1.76%  "

$ws.Range("D16").Value = "'3.795.24"
$ws.Range("E16").Value = "'  -4.73%  "

$ws.Range("D17").Value = "'65.556.73"
$ws.Range("E17").Value = "'  -4.51%  "

$ws.Range("D18").Value = "'17.81"
$ws.Range("E18").Value = "'  -0.59%  "

$ws.Range("E19").Value = "'  -3.21%  "

$ws.Range("D20").Value = "'3.269.23"
$ws.Range("E20").Value = "'  -4.64%  "

$ws.Range("E21").Value = "'  -6.38%  "

$ws.Range("D22").Value = "'0.902"
$ws.Range("E22").Value = "'  -3.17%  "

$ws.Range("D23").Value = "'17.59"
$ws.Range("E23").Value = "'  +0.16%  "

$ws.Range("D24").Value = "'106.13"
$ws.Range("E24").Value = "'  +8.73%  "

$ws.Range("E25").Value = "'  -6.58%  "

$ws.Range("E26").Value = "'  -6.32%  "

$ws.Range("E27").Value = "'  -5.62%  "

$ws.Range("D28").Value = "'9.53"
$ws.Range("E28").Value = "'  -2.49%  "

$ws.Range("D29").Value = "'8.63"
$ws.Range("E29").Value = "'  -5.29%  "

$ws.Range("D30").Value = "'30.31"
$ws.Range("E30").Value = "'  -5.23%  "

$ws.Range("D31").Value = "'4.04"
$ws.Range("E31").Value = "'  -3.59%  "

$ws.Range("D32").Value = "'6.24"
$ws.Range("E32").Value = "'  -5.58%  "

$ws.Range("D33").Value = "'10.98"
$ws.Range("E33").Value = "'  -3.97%  "

$ws.Range("D34").Value = "'540.85"
$ws.Range("E34").Value = "'  +10.84%  "

$ws.Range("D35").Value = "'0.104"
$ws.Range("E35").Value = "'  -2.94%  "

$ws.Range("E36").Value = "'  -0.21%  "

$ws.Range("D37").Value = "'57.21"
$ws.Range("E37").Value = "'  -5.84%  "

$ws.Range("D38").Value = "'3.666.44"
$ws.Range("E38").Value = "'  +0.83%  "

$ws.Range("D39").Value = "'3.37"
$ws.Range("E39").Value = "'  -1.03%  "

$ws.Range("D40").Value = "'0.0₃0719"
$ws.Range("E40").Value = "'  -7.20%  "

$ws.Range("E41").Value = "'  -1.21%  "

$ws.Range("E42").Value = "'  -5.28%  "

$ws.Range("B43").Value = "'InjectiveProtocol"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'32.53"
$ws.Range("E43").Value = "'  -3.77%  "

$ws.Range("B44").Value = "'CoreDAO"
$ws.Range("C44").Value = "'https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").Value = "'3.32"
$ws.Range("E44").Value = "'  -5.87%  "

$ws.Range("D45").Value = "'0.335"
$ws.Range("E45").Value = "'  -7.83%  "

$ws.Range("B46").Value = "'ApeXProtocol"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.24"
$ws.Range("E46").Value = "'  -2.16%  "

$ws.Range("B47").Value = "'VeChain"
$ws.Range("C47").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0413"
$ws.Range("E47").Value = "'  -4.57%  "

$ws.Range("E48").Value = "'  -5.89%  "

$ws.Range("E49").Value = "'  -3.16%  "

$ws.Range("E50").Value = "'  -0.04%  "

$ws.Range("E51").Value = "'  +1.90%  "
